$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Campus" column (column B); remaining columns shift left
# (C->B, D->C, ..., J->I), matching the source-data refresh in the diff.
$ws.Columns("B:B").Delete()

# Refresh the cells whose source data changed (using the new, post-shift
# column letters).
$ws.Range('G2').Value = '42,360 (fall 2025)'
$ws.Range('C3').Value = 'March 27, 1871; 154 years ago (1871-03-27)'
$ws.Range('D3').Value = 'Fayetteville, Arkansas, United States36°4′7″N 94°10′34″W﻿ / ﻿36.06861°N 94.17611°W﻿ / 36.06861; -94.17611'
$ws.Range('G3').Value = '34,174 (fall 2025)'
$ws.Range('B4').Value = '$1.25 billion (FY2024)'
$ws.Range('B6').Value = '$2.056 billion (2024)'
$ws.Range('D6').Value = 'Athens, Georgia, US33°57′21″N 83°22′28″W﻿ / ﻿33.9558°N 83.3745°W﻿ / 33.9558; -83.3745'
$ws.Range('B7').Value = '$2.27 billion (2024)'
$ws.Range('C7').Value = 'February 22, 1865;160 years ago (1865-02-22)'
$ws.Range('G8').Value = '42,016 (fall 2024)'
$ws.Range('B9').Value = '$1 billion (2025)'
$ws.Range('D9').Value = 'Oxford, Mississippi'
$ws.Range('B10').Value = '$1.0 billion (2024)'
$ws.Range('F12').Value = 'NCAA Division I FBS – SECBig 12MPSF'
$ws.Range('C14').Value = 'September 10, 1794; 231 years ago (1794-09-10)'
$ws.Range('G14').Value = '40,784 (fall 2025)'
$ws.Range('B15').Value = '$20.85 billion (FY2024)(UT Austin only)$47.47 billion (FY2024)(system-wide)'
$ws.Range('C15').Value = 'September 15, 1883; 142 years ago (1883-09-15)'
$ws.Range('G15').Value = '53,864 (fall 2024)'
$ws.Range('G16').Value = '79,114 (fall 2024) • 71,045 (College Station) • 2,138 (Galveston) • 1,751 (Fort Worth) • 430 (McAllen) • 3,750 (Health Science Center)'
$ws.Range('B17').Value = '$10.2 billion (2024)'

# "27,124" would otherwise be auto-parsed as the number 27124 by Excel
# (it has no non-numeric text to keep it a string), so mark it as text
# first to preserve it verbatim, matching the source data.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '27,124'
